# fix heights of title text boxes to match wrapped text
#
# Task Id=1 ("line wrap long descriptions and increase textbox height to
# accomodate") has been completed, so it moves from the "Active" sheet
# (status Todo) to the "Inactive" sheet as a "Done" item, stamped with the
# date it was completed.

$wb  = $excel.ActiveWorkbook
$active   = $wb.Worksheets.Item("Active")
$inactive = $wb.Worksheets.Item("Inactive")

# Make room for the completed task at the top of the Inactive list (it keeps
# the lowest Id, so it sorts above the existing "program icon " row).
$inactive.Rows.Item(2).Insert()

# Copy the whole row across sheets so every cell keeps its original raw
# type/value (in particular the "Created" date stays plain text instead of
# being re-interpreted as a serial date number).
$active.Range("A2:F2").Copy($inactive.Range("A2:F2"))

# It is now Done rather than Todo.
$inactive.Range("C2").Value = "Done"

# Stamp the date it was marked Done. Force the cell to Text first so Excel
# doesn't reinterpret the typed string as a date serial, then drop back to
# the Normal style so no explicit number format sticks to the cell.
$inactive.Range("F2").NumberFormat = "@"
$inactive.Range("F2").Value = "3/2/2018"
$inactive.Range("F2").Style = "Normal"

# Remove the now-completed task from the Active sheet.
$active.Rows.Item(2).Delete()
